$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "37.047.04"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.012.78"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -1.97%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "225.63"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.604"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("E7").Value = "  -0.01%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "54.88"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -4.10%  "

$ws.Range("E9").Value = "  -3.43%  "

$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("E11").Value = "  -5.13%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "2.308.54"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "14.02"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -4.95%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "19.79"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -5.15%  "

$ws.Range("E15").Value = "  -2.39%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.736"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -3.41%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "2.021.09"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "37.006.70"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "6.28"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +2.96%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "68.23"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0813"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -3.50%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "222.44"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("E24").Value = "  +1.32%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.17"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -5.58%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "164.44"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -2.04%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "8.99"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -6.30%  "

$ws.Range("E28").Value = "  -4.32%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "18.55"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -2.51%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "1.30"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -8.46%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.117"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -1.69%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.46"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("E33").Value = "  -2.99%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0600"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -2.80%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "2.32"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").Value = "  -4.71%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "5.35"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -0.80%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "1.455.29"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -2.52%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0212"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -4.69%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "94.86"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -2.10%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "2.76"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -4.55%  "

$ws.Range("E44").Value = "  -4.58%  "

$ws.Range("E45").Value = "  -4.62%  "

$ws.Range("E46").Value = "  -7.95%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "7.11"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "0.998"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -2.73%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "2.90"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "2.193.53"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -2.27%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "3.59"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -3.15%  "
